# Insert a new weekly price record as row 83 ("Hortaliza, Terminal
# Hortofrutícola Agro Chillán - Arveja Verde"), pushing the existing
# rows 83:102 down to 84:103 (dimension grows from A1:R102 to A1:R103).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 83:102 down one row, leaving a blank row 83 for the new record.
$ws.Rows.Item(83).Insert()

# Populate the new row 83 with the latest weekly observation.
$ws.Cells.Item(83, 1).Value  = 7
$ws.Cells.Item(83, 2).Value  = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(83, 3).Value  = "Ñuble"
$ws.Cells.Item(83, 4).Value  = 44943
$ws.Cells.Item(83, 5).Value  = 16
$ws.Cells.Item(83, 6).Value  = 100112022
$ws.Cells.Item(83, 7).Value  = "Arveja Verde"
$ws.Cells.Item(83, 8).Value  = "Sin especificar"
$ws.Cells.Item(83, 9).Value  = "Primera"
$ws.Cells.Item(83, 10).Value = 60
$ws.Cells.Item(83, 11).Value = 22000
$ws.Cells.Item(83, 12).Value = 23000
$ws.Cells.Item(83, 13).Value = 22500
$ws.Cells.Item(83, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(83, 15).Value = "Región de Ñuble"
$ws.Cells.Item(83, 16).Value = 900
$ws.Cells.Item(83, 17).Value = 25
$ws.Cells.Item(83, 18).Value = "Hortaliza"
